$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before X to make room for the Year column
$ws.Columns.Item(24).Insert()

# Set the new header
$ws.Range("X1").Value = "Year"

# Populate the Year values for each movie row
$ws.Range("X2").Value = 2015
$ws.Range("X3").Value = 1974
$ws.Range("X4").Value = 2005
$ws.Range("X5").Value = 2004
$ws.Range("X6").Value = 2011
$ws.Range("X7").Value = 2002
$ws.Range("X8").Value = 2017
$ws.Range("X9").Value = 2013
$ws.Range("X10").Value = 2015
$ws.Range("X11").Value = 2021
$ws.Range("X12").Value = 2016
$ws.Range("X13").Value = 2016
$ws.Range("X14").Value = 2000
$ws.Range("X15").Value = 1986
$ws.Range("X16").Value = 2018
$ws.Range("X17").Value = 2016
$ws.Range("X18").Value = 2008
$ws.Range("X19").Value = 2005
$ws.Range("X20").Value = 2015
$ws.Range("X21").Value = 1975
$ws.Range("X22").Value = 1991
$ws.Range("X23").Value = 2009
$ws.Range("X24").Value = 2002
$ws.Range("X25").Value = 2004
$ws.Range("X26").Value = 1971
$ws.Range("X27").Value = 2006
$ws.Range("X28").Value = 2018
$ws.Range("X29").Value = 1939
$ws.Range("X30").Value = 1996
$ws.Range("X31").Value = 2020
$ws.Range("X32").Value = 2014
$ws.Range("X33").Value = 1984
$ws.Range("X34").Value = 2017
$ws.Range("X35").Value = 2015
$ws.Range("X36").Value = 1997
$ws.Range("X37").Value = 1986
$ws.Range("X38").Value = 2019
$ws.Range("X39").Value = 2009
$ws.Range("X40").Value = 2012
$ws.Range("X41").Value = 1991
$ws.Range("X42").Value = 2016
$ws.Range("X43").Value = 2000
$ws.Range("X44").Value = 2007
$ws.Range("X45").Value = 2016
$ws.Range("X46").Value = 2017
$ws.Range("X47").Value = 1999
$ws.Range("X48").Value = 1992
$ws.Range("X49").Value = 1998
$ws.Range("X50").Value = 1987
$ws.Range("X51").Value = 2016
$ws.Range("X52").Value = 2008
$ws.Range("X53").Value = 2008
$ws.Range("X54").Value = 1988
$ws.Range("X55").Value = 2015
$ws.Range("X56").Value = 2017
$ws.Range("X57").Value = 2013
$ws.Range("X58").Value = 2008
$ws.Range("X59").Value = 2018
$ws.Range("X60").Value = 2012
$ws.Range("X61").Value = 2012
$ws.Range("X62").Value = 2011
$ws.Range("X63").Value = 2017
$ws.Range("X64").Value = 2018
$ws.Range("X65").Value = 2012
$ws.Range("X66").Value = 2006
$ws.Range("X67").Value = 2010
$ws.Range("X68").Value = 2019
$ws.Range("X69").Value = 2018
$ws.Range("X70").Value = 2015
$ws.Range("X71").Value = 2020
$ws.Range("X72").Value = 2011
$ws.Range("X73").Value = 1982
$ws.Range("X74").Value = 2009
$ws.Range("X75").Value = 2019
$ws.Range("X76").Value = 1985
$ws.Range("X77").Value = 2012
$ws.Range("X78").Value = 2010
$ws.Range("X79").Value = 2019
$ws.Range("X80").Value = 2013
$ws.Range("X81").Value = 2005
$ws.Range("X82").Value = 2002
$ws.Range("X83").Value = 1996
$ws.Range("X84").Value = 2003
$ws.Range("X85").Value = 2006
$ws.Range("X86").Value = 2021
$ws.Range("X87").Value = 1989
$ws.Range("X88").Value = 2017
$ws.Range("X89").Value = 2015
$ws.Range("X90").Value = 1986
$ws.Range("X91").Value = 1987
$ws.Range("X92").Value = 2019
$ws.Range("X93").Value = 1975
$ws.Range("X94").Value = 2011
$ws.Range("X95").Value = 2004
$ws.Range("X96").Value = 1980
$ws.Range("X97").Value = 2004
$ws.Range("X98").Value = 2018
$ws.Range("X99").Value = 1984
$ws.Range("X100").Value = 1976
$ws.Range("X101").Value = 2009
$ws.Range("X102").Value = 2010
$ws.Range("X103").Value = 1996
$ws.Range("X104").Value = 2006
$ws.Range("X105").Value = 2009
$ws.Range("X106").Value = 2006
$ws.Range("X107").Value = 2009
$ws.Range("X108").Value = 2007
$ws.Range("X109").Value = 2014
$ws.Range("X110").Value = 1999
$ws.Range("X111").Value = 1983
$ws.Range("X112").Value = 1997
$ws.Range("X113").Value = 2012
$ws.Range("X114").Value = 1999
$ws.Range("X115").Value = 2018
$ws.Range("X116").Value = 1942
$ws.Range("X117").Value = 1986
$ws.Range("X118").Value = 1985
$ws.Range("X119").Value = 1986
$ws.Range("X120").Value = 1992
$ws.Range("X121").Value = 2014
$ws.Range("X122").Value = 2008
$ws.Range("X123").Value = 1993
$ws.Range("X124").Value = 1968
$ws.Range("X125").Value = 1995
$ws.Range("X126").Value = 2006
$ws.Range("X127").Value = 2017
$ws.Range("X128").Value = 1991
$ws.Range("X129").Value = 2005
$ws.Range("X130").Value = 1996
$ws.Range("X131").Value = 2019
$ws.Range("X132").Value = 2007
$ws.Range("X133").Value = 2007
$ws.Range("X134").Value = 2017
